$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 5899.727  # H40: 5616.4165 -> 5899.727
$ws.Cells.Item(40, 10).Value = 2483  # J40: 2485.4285 -> 2483
$ws.Cells.Item(40, 12).Value = 2483  # L40: 2485.4285 -> 2483
$ws.Cells.Item(40, 14).Value = -2833  # N40: -2835.4285 -> -2833
$ws.Cells.Item(111, 8).Value = 1619  # H111: 1730.25 -> 1619
$ws.Cells.Item(111, 9).Value = 1023.75  # I111: 1590.6666 -> 1023.75
$ws.Cells.Item(111, 10).Value = 4000  # J111: 2149 -> 4000
$ws.Cells.Item(111, 11).Value = 3071.25  # K111: 4771.9998 -> 3071.25
$ws.Cells.Item(111, 12).Value = 12000  # L111: 6447 -> 12000
$ws.Cells.Item(111, 13).Value = -4.25  # M111: -1704.9998 -> -4.25
$ws.Cells.Item(111, 14).Value = -18134  # N111: -12581 -> -18134
$ws.Cells.Item(112, 8).Value = 2571.1667  # H112: 2573.9443 -> 2571.1667
$ws.Cells.Item(112, 10).Value = 2848.3103  # J112: 2851.7585 -> 2848.3103
$ws.Cells.Item(112, 12).Value = 8544.930899999999  # L112: 8555.2755 -> 8544.930899999999
$ws.Cells.Item(112, 14).Value = -10760.9309  # N112: -10771.2755 -> -10760.9309
$ws.Cells.Item(137, 8).Value = 3037.68  # H137: 3143.4583 -> 3037.68
$ws.Cells.Item(137, 9).Value = 1728.2  # I137: 1864.7778 -> 1728.2
$ws.Cells.Item(137, 11).Value = 5184.6  # K137: 5594.3334 -> 5184.6
$ws.Cells.Item(137, 13).Value = -2634.6  # M137: -3044.3334 -> -2634.6
$ws.Cells.Item(138, 8).Value = 2939.9246  # H138: 2981.7646 -> 2939.9246
$ws.Cells.Item(138, 10).Value = 2979.1365  # J138: 3031.8096 -> 2979.1365
$ws.Cells.Item(138, 12).Value = 8937.4095  # L138: 9095.4288 -> 8937.4095
$ws.Cells.Item(138, 14).Value = -19217.4095  # N138: -19375.4288 -> -19217.4095

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 899.9245  # H2: 899.9434 -> 899.9245
$ws.Cells.Item(2, 9).Value = 827.1667  # I2: 827.1905 -> 827.1667
$ws.Cells.Item(2, 11).Value = 827.1667  # K2: 827.1905 -> 827.1667
$ws.Cells.Item(2, 13).Value = -714.1667  # M2: -714.1905 -> -714.1667
$ws.Cells.Item(45, 8).Value = 922.9  # H45: 923.5 -> 922.9
$ws.Cells.Item(45, 9).Value = 922.9  # I45: 923.5 -> 922.9
$ws.Cells.Item(45, 11).Value = 922.9  # K45: 923.5 -> 922.9
$ws.Cells.Item(45, 13).Value = -545.9  # M45: -546.5 -> -545.9
$ws.Cells.Item(61, 8).Value = 3802.7  # H61: 3452.238 -> 3802.7
$ws.Cells.Item(61, 9).Value = 3513  # I61: 3206.4443 -> 3513
$ws.Cells.Item(61, 10).Value = 4961.5  # J61: 4927 -> 4961.5
$ws.Cells.Item(61, 11).Value = 3513  # K61: 3206.4443 -> 3513
$ws.Cells.Item(61, 12).Value = 4961.5  # L61: 4927 -> 4961.5
$ws.Cells.Item(61, 13).Value = -3301  # M61: -2994.4443 -> -3301
$ws.Cells.Item(61, 14).Value = -5385.5  # N61: -5351 -> -5385.5
$ws.Cells.Item(63, 8).Value = 2489.2942  # H63: 2519.08 -> 2489.2942
$ws.Cells.Item(63, 10).Value = 2555.375  # J63: 2777.5715 -> 2555.375
$ws.Cells.Item(63, 12).Value = 2555.375  # L63: 2777.5715 -> 2555.375
$ws.Cells.Item(63, 14).Value = -3927.375  # N63: -4149.5715 -> -3927.375
$ws.Cells.Item(66, 8).Value = 2489.2942  # H66: 2519.08 -> 2489.2942
$ws.Cells.Item(66, 10).Value = 2555.375  # J66: 2777.5715 -> 2555.375
$ws.Cells.Item(66, 12).Value = 12776.875  # L66: 13887.8575 -> 12776.875
$ws.Cells.Item(66, 14).Value = -19640.875  # N66: -20751.8575 -> -19640.875
$ws.Cells.Item(102, 8).Value = 2264.4348  # H102: 2272.8696 -> 2264.4348
$ws.Cells.Item(102, 9).Value = 2255.25  # I102: 2358.5264 -> 2255.25
$ws.Cells.Item(102, 10).Value = 2325.6667  # J102: 1866 -> 2325.6667
$ws.Cells.Item(102, 11).Value = 2255.25  # K102: 2358.5264 -> 2255.25
$ws.Cells.Item(102, 12).Value = 2325.6667  # L102: 1866 -> 2325.6667
$ws.Cells.Item(102, 13).Value = -633.25  # M102: -736.5264000000002 -> -633.25
$ws.Cells.Item(102, 14).Value = -5569.6667  # N102: -5110 -> -5569.6667
$ws.Cells.Item(116, 8).Value = 899.9245  # H116: 899.9434 -> 899.9245
$ws.Cells.Item(116, 9).Value = 827.1667  # I116: 827.1905 -> 827.1667
$ws.Cells.Item(116, 11).Value = 827.1667  # K116: 827.1905 -> 827.1667
$ws.Cells.Item(116, 13).Value = 1466.8333  # M116: 1466.8095 -> 1466.8333
$ws.Cells.Item(132, 8).Value = 1910  # H132: 1521.6666 -> 1910
$ws.Cells.Item(132, 9).Value = 1910  # I132: 1521.6666 -> 1910
$ws.Cells.Item(132, 11).Value = 5730  # K132: 4564.9998 -> 5730
$ws.Cells.Item(132, 13).Value = -3200  # M132: -2034.9998 -> -3200
$ws.Cells.Item(136, 8).Value = 3802.7  # H136: 3452.238 -> 3802.7
$ws.Cells.Item(136, 9).Value = 3513  # I136: 3206.4443 -> 3513
$ws.Cells.Item(136, 10).Value = 4961.5  # J136: 4927 -> 4961.5
$ws.Cells.Item(136, 11).Value = 10539  # K136: 9619.332900000001 -> 10539
$ws.Cells.Item(136, 12).Value = 14884.5  # L136: 14781 -> 14884.5
$ws.Cells.Item(136, 13).Value = -7989  # M136: -7069.332900000001 -> -7989
$ws.Cells.Item(136, 14).Value = -19984.5  # N136: -19881 -> -19984.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 899.9245  # H3: 899.9434 -> 899.9245
$ws.Cells.Item(3, 9).Value = 827.1667  # I3: 827.1905 -> 827.1667
$ws.Cells.Item(3, 11).Value = 827.1667  # K3: 827.1905 -> 827.1667
$ws.Cells.Item(3, 13).Value = -713.1667  # M3: -713.1905 -> -713.1667
$ws.Cells.Item(99, 8).Value = 1215.9412  # H99: 1239 -> 1215.9412
$ws.Cells.Item(99, 9).Value = 1029.3077  # I99: 1065.3334 -> 1029.3077
$ws.Cells.Item(99, 10).Value = 1822.5  # J99: 1655.8 -> 1822.5
$ws.Cells.Item(99, 11).Value = 1029.3077  # K99: 1065.3334 -> 1029.3077
$ws.Cells.Item(99, 12).Value = 1822.5  # L99: 1655.8 -> 1822.5
$ws.Cells.Item(99, 13).Value = 468.6922999999999  # M99: 432.6666 -> 468.6922999999999
$ws.Cells.Item(99, 14).Value = -4818.5  # N99: -4651.8 -> -4818.5
$ws.Cells.Item(105, 8).Value = 2085  # H105: 2112.8572 -> 2085
$ws.Cells.Item(105, 10).Value = 2212.25  # J105: 2314 -> 2212.25
$ws.Cells.Item(105, 12).Value = 2212.25  # L105: 2314 -> 2212.25
$ws.Cells.Item(105, 14).Value = -5706.25  # N105: -5808 -> -5706.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 859  # H16: 916.9 -> 859
$ws.Cells.Item(16, 9).Value = 816.2353000000001  # I16: 871.26666 -> 816.2353000000001
$ws.Cells.Item(16, 10).Value = 980.1667  # J16: 1053.8 -> 980.1667
$ws.Cells.Item(16, 11).Value = 816.2353000000001  # K16: 871.26666 -> 816.2353000000001
$ws.Cells.Item(16, 12).Value = 980.1667  # L16: 1053.8 -> 980.1667
$ws.Cells.Item(16, 13).Value = -529.2353000000001  # M16: -584.26666 -> -529.2353000000001
$ws.Cells.Item(16, 14).Value = -1554.1667  # N16: -1627.8 -> -1554.1667
$ws.Cells.Item(22, 8).Value = 495  # H22: 496.66666 -> 495
$ws.Cells.Item(22, 9).Value = 500  # I22: 496.66666 -> 500
$ws.Cells.Item(22, 10).Value = 490  # J22: 0 -> 490
$ws.Cells.Item(22, 11).Value = 500  # K22: 496.66666 -> 500
$ws.Cells.Item(22, 12).Value = 490  # L22: 0 -> 490
$ws.Cells.Item(22, 13).Value = -150  # M22: -146.66666 -> -150
$ws.Cells.Item(31, 8).Value = 5102.4355  # H31: 5035.127 -> 5102.4355
$ws.Cells.Item(31, 9).Value = 3387.1428  # I31: 3240.8647 -> 3387.1428
$ws.Cells.Item(31, 10).Value = 7325.963  # J31: 7588.5 -> 7325.963
$ws.Cells.Item(31, 11).Value = 3387.1428  # K31: 3240.8647 -> 3387.1428
$ws.Cells.Item(31, 12).Value = 7325.963  # L31: 7588.5 -> 7325.963
$ws.Cells.Item(31, 13).Value = -3092.1428  # M31: -2945.8647 -> -3092.1428
$ws.Cells.Item(31, 14).Value = -7915.963  # N31: -8178.5 -> -7915.963
$ws.Cells.Item(34, 8).Value = 5102.4355  # H34: 5035.127 -> 5102.4355
$ws.Cells.Item(34, 9).Value = 3387.1428  # I34: 3240.8647 -> 3387.1428
$ws.Cells.Item(34, 10).Value = 7325.963  # J34: 7588.5 -> 7325.963
$ws.Cells.Item(34, 11).Value = 3387.1428  # K34: 3240.8647 -> 3387.1428
$ws.Cells.Item(34, 12).Value = 7325.963  # L34: 7588.5 -> 7325.963
$ws.Cells.Item(34, 13).Value = -3185.1428  # M34: -3038.8647 -> -3185.1428
$ws.Cells.Item(34, 14).Value = -7729.963  # N34: -7992.5 -> -7729.963
$ws.Cells.Item(55, 8).Value = 0  # H55: 20073 -> 0
$ws.Cells.Item(55, 9).Value = 0  # I55: 20073 -> 0
$ws.Cells.Item(55, 11).Value = 0  # K55: 20073 -> 0
$ws.Cells.Item(58, 8).Value = 4523.385  # H58: 5163.909 -> 4523.385
$ws.Cells.Item(58, 9).Value = 4523.385  # I58: 5163.909 -> 4523.385
$ws.Cells.Item(58, 11).Value = 4523.385  # K58: 5163.909 -> 4523.385
$ws.Cells.Item(58, 13).Value = -4320.385  # M58: -4960.909 -> -4320.385
$ws.Cells.Item(86, 8).Value = 7853.5  # H86: 6443 -> 7853.5
$ws.Cells.Item(86, 9).Value = 0  # I86: 4562.3335 -> 0
$ws.Cells.Item(86, 11).Value = 0  # K86: 4562.3335 -> 0
$ws.Cells.Item(89, 8).Value = 7853.5  # H89: 6443 -> 7853.5
$ws.Cells.Item(89, 9).Value = 0  # I89: 4562.3335 -> 0
$ws.Cells.Item(89, 11).Value = 0  # K89: 22811.6675 -> 0
$ws.Cells.Item(107, 8).Value = 2899285.8  # H107: 2299476 -> 2899285.8
$ws.Cells.Item(107, 9).Value = 4167361.5  # I107: 3175184.2 -> 4167361.5
$ws.Cells.Item(107, 10).Value = 826.5714  # J107: 742 -> 826.5714
$ws.Cells.Item(107, 11).Value = 4167361.5  # K107: 3175184.2 -> 4167361.5
$ws.Cells.Item(107, 12).Value = 826.5714  # L107: 742 -> 826.5714
$ws.Cells.Item(107, 13).Value = -4165441.5  # M107: -3173264.2 -> -4165441.5
$ws.Cells.Item(107, 14).Value = -4666.5714  # N107: -4582 -> -4666.5714
$ws.Cells.Item(113, 8).Value = 859  # H113: 916.9 -> 859
$ws.Cells.Item(113, 9).Value = 816.2353000000001  # I113: 871.26666 -> 816.2353000000001
$ws.Cells.Item(113, 10).Value = 980.1667  # J113: 1053.8 -> 980.1667
$ws.Cells.Item(113, 11).Value = 816.2353000000001  # K113: 871.26666 -> 816.2353000000001
$ws.Cells.Item(113, 12).Value = 980.1667  # L113: 1053.8 -> 980.1667
$ws.Cells.Item(113, 13).Value = 1353.7647  # M113: 1298.73334 -> 1353.7647
$ws.Cells.Item(113, 14).Value = -5320.1667  # N113: -5393.8 -> -5320.1667
$ws.Cells.Item(132, 8).Value = 2189.1365  # H132: 2561.3809 -> 2189.1365
$ws.Cells.Item(132, 9).Value = 2245.0952  # I132: 2464.5 -> 2245.0952
$ws.Cells.Item(132, 10).Value = 1014  # J132: 4499 -> 1014
$ws.Cells.Item(132, 11).Value = 6735.285600000001  # K132: 7393.5 -> 6735.285600000001
$ws.Cells.Item(132, 12).Value = 3042  # L132: 13497 -> 3042
$ws.Cells.Item(132, 13).Value = -4205.285600000001  # M132: -4863.5 -> -4205.285600000001
$ws.Cells.Item(132, 14).Value = -8102  # N132: -18557 -> -8102
$ws.Cells.Item(134, 8).Value = 1297.9667  # H134: 1307.6774 -> 1297.9667
$ws.Cells.Item(134, 9).Value = 1307.875  # I134: 1319.52 -> 1307.875
$ws.Cells.Item(134, 11).Value = 3923.625  # K134: 3958.56 -> 3923.625
$ws.Cells.Item(134, 13).Value = -1388.625  # M134: -1423.56 -> -1388.625
$ws.Cells.Item(136, 8).Value = 4523.385  # H136: 5163.909 -> 4523.385
$ws.Cells.Item(136, 9).Value = 4523.385  # I136: 5163.909 -> 4523.385
$ws.Cells.Item(136, 11).Value = 13570.155  # K136: 15491.727 -> 13570.155
$ws.Cells.Item(136, 13).Value = -11020.155  # M136: -12941.727 -> -11020.155
$ws.Cells.Item(22, 14).Value = -1190  # N22: (new) -> -1190
$ws.Cells.Item(55, 13).ClearContents()  # M55: -19758 -> (cleared)
$ws.Cells.Item(86, 13).ClearContents()  # M86: -3439.3335 -> (cleared)
$ws.Cells.Item(89, 13).ClearContents()  # M89: -17195.6675 -> (cleared)

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(75, 8).Value = 3904.1667  # H75: 3904.4167 -> 3904.1667
$ws.Cells.Item(75, 10).Value = 4426.6665  # J75: 4427 -> 4426.6665
$ws.Cells.Item(75, 12).Value = 13279.9995  # L75: 13281 -> 13279.9995
$ws.Cells.Item(75, 14).Value = -15275.9995  # N75: -15277 -> -15275.9995
$ws.Cells.Item(78, 8).Value = 3904.1667  # H78: 3904.4167 -> 3904.1667
$ws.Cells.Item(78, 10).Value = 4426.6665  # J78: 4427 -> 4426.6665
$ws.Cells.Item(78, 12).Value = 39839.9985  # L78: 39843 -> 39839.9985
$ws.Cells.Item(78, 14).Value = -49823.9985  # N78: -49827 -> -49823.9985
$ws.Cells.Item(107, 8).Value = 564.1053000000001  # H107: 533.6667 -> 564.1053000000001
$ws.Cells.Item(107, 9).Value = 306.625  # I107: 294.66666 -> 306.625
$ws.Cells.Item(107, 10).Value = 751.36365  # J107: 712.9167 -> 751.36365
$ws.Cells.Item(107, 11).Value = 919.875  # K107: 883.9999799999999 -> 919.875
$ws.Cells.Item(107, 12).Value = 2254.09095  # L107: 2138.7501 -> 2254.09095
$ws.Cells.Item(107, 13).Value = 1000.125  # M107: 1036.00002 -> 1000.125
$ws.Cells.Item(107, 14).Value = -6094.09095  # N107: -5978.7501 -> -6094.09095
$ws.Cells.Item(122, 8).Value = 1543.48  # H122: 1524.375 -> 1543.48
$ws.Cells.Item(122, 10).Value = 1600.4117  # J122: 1575.3125 -> 1600.4117
$ws.Cells.Item(122, 12).Value = 14403.7053  # L122: 14177.8125 -> 14403.7053
$ws.Cells.Item(122, 14).Value = -19303.7053  # N122: -19077.8125 -> -19303.7053
$ws.Cells.Item(137, 8).Value = 7171.8184  # H137: 6894.5 -> 7171.8184
$ws.Cells.Item(137, 10).Value = 2840  # J137: 3007.3333 -> 2840
$ws.Cells.Item(137, 12).Value = 8520  # L137: 9021.999899999999 -> 8520
$ws.Cells.Item(137, 14).Value = -18720  # N137: -19221.9999 -> -18720
$ws.Cells.Item(139, 8).Value = 9337.272000000001  # H139: 7448.619 -> 9337.272000000001
$ws.Cells.Item(139, 9).Value = 8067.5386  # I139: 4656.5835 -> 8067.5386
$ws.Cells.Item(139, 11).Value = 24202.6158  # K139: 13969.7505 -> 24202.6158
$ws.Cells.Item(139, 13).Value = -19062.6158  # M139: -8829.750499999998 -> -19062.6158
$ws.Cells.Item(140, 8).Value = 2054.7812  # H140: 2098.484 -> 2054.7812
$ws.Cells.Item(140, 9).Value = 1219.5385  # I140: 1262.8334 -> 1219.5385
$ws.Cells.Item(140, 11).Value = 3658.6155  # K140: 3788.5002 -> 3658.6155
$ws.Cells.Item(140, 13).Value = 1521.3845  # M140: 1391.4998 -> 1521.3845

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 2119.8823  # H113: 2164.8125 -> 2119.8823
$ws.Cells.Item(113, 9).Value = 2082.5334  # I113: 2109.1333 -> 2082.5334
$ws.Cells.Item(113, 10).Value = 2400  # J113: 3000 -> 2400
$ws.Cells.Item(113, 11).Value = 2082.5334  # K113: 2109.1333 -> 2082.5334
$ws.Cells.Item(113, 12).Value = 2400  # L113: 3000 -> 2400
$ws.Cells.Item(113, 13).Value = 87.4666000000002  # M113: 60.86670000000004 -> 87.4666000000002
$ws.Cells.Item(113, 14).Value = -6740  # N113: -7340 -> -6740
$ws.Cells.Item(122, 8).Value = 68444.13  # H122: 60608.824 -> 68444.13
$ws.Cells.Item(122, 9).Value = 168116  # I122: 144342.28 -> 168116
$ws.Cells.Item(122, 10).Value = 1996.2222  # J122: 1995.4 -> 1996.2222
$ws.Cells.Item(122, 11).Value = 504348  # K122: 433026.84 -> 504348
$ws.Cells.Item(122, 12).Value = 5988.6666  # L122: 5986.200000000001 -> 5988.6666
$ws.Cells.Item(122, 13).Value = -501898  # M122: -430576.84 -> -501898
$ws.Cells.Item(122, 14).Value = -10888.6666  # N122: -10886.2 -> -10888.6666
$ws.Cells.Item(132, 8).Value = 1755.1333  # H132: 1785.9318 -> 1755.1333
$ws.Cells.Item(132, 9).Value = 1542.3143  # I132: 1575.9117 -> 1542.3143
$ws.Cells.Item(132, 11).Value = 4626.9429  # K132: 4727.7351 -> 4626.9429
$ws.Cells.Item(132, 13).Value = -2096.9429  # M132: -2197.7351 -> -2096.9429
$ws.Cells.Item(136, 8).Value = 13539.728  # H136: 13480.25 -> 13539.728
$ws.Cells.Item(136, 10).Value = 13539.728  # J136: 13480.25 -> 13539.728
$ws.Cells.Item(136, 12).Value = 40619.18399999999  # L136: 40440.75 -> 40619.18399999999
$ws.Cells.Item(136, 14).Value = -45719.18399999999  # N136: -45540.75 -> -45719.18399999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 9019.385  # H7: 7661.294 -> 9019.385
$ws.Cells.Item(7, 9).Value = 8123.8184  # I7: 7119 -> 8123.8184
$ws.Cells.Item(7, 10).Value = 13945  # J7: 9423.75 -> 13945
$ws.Cells.Item(7, 11).Value = 8123.8184  # K7: 7119 -> 8123.8184
$ws.Cells.Item(7, 12).Value = 13945  # L7: 9423.75 -> 13945
$ws.Cells.Item(7, 13).Value = -8011.8184  # M7: -7007 -> -8011.8184
$ws.Cells.Item(7, 14).Value = -14169  # N7: -9647.75 -> -14169
$ws.Cells.Item(46, 8).Value = 1476.5834  # H46: 1389.1786 -> 1476.5834
$ws.Cells.Item(46, 9).Value = 2248.1667  # I46: 1869.75 -> 2248.1667
$ws.Cells.Item(46, 10).Value = 1219.3889  # J46: 1196.95 -> 1219.3889
$ws.Cells.Item(46, 11).Value = 2248.1667  # K46: 1869.75 -> 2248.1667
$ws.Cells.Item(46, 12).Value = 1219.3889  # L46: 1196.95 -> 1219.3889
$ws.Cells.Item(46, 13).Value = -2060.1667  # M46: -1681.75 -> -2060.1667
$ws.Cells.Item(46, 14).Value = -1595.3889  # N46: -1572.95 -> -1595.3889
$ws.Cells.Item(108, 8).Value = 89000.5  # H108: 89001 -> 89000.5
$ws.Cells.Item(108, 10).Value = 89000.5  # J108: 89001 -> 89000.5
$ws.Cells.Item(108, 12).Value = 89000.5  # L108: 89001 -> 89000.5
$ws.Cells.Item(108, 14).Value = -96680.5  # N108: -96681 -> -96680.5
$ws.Cells.Item(126, 8).Value = 9019.385  # H126: 7661.294 -> 9019.385
$ws.Cells.Item(126, 9).Value = 8123.8184  # I126: 7119 -> 8123.8184
$ws.Cells.Item(126, 10).Value = 13945  # J126: 9423.75 -> 13945
$ws.Cells.Item(126, 11).Value = 24371.4552  # K126: 21357 -> 24371.4552
$ws.Cells.Item(126, 12).Value = 41835  # L126: 28271.25 -> 41835
$ws.Cells.Item(126, 13).Value = -21901.4552  # M126: -18887 -> -21901.4552
$ws.Cells.Item(126, 14).Value = -46775  # N126: -33211.25 -> -46775

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 16052.8  # H41: 18746.75 -> 16052.8
$ws.Cells.Item(41, 10).Value = 24137  # J41: 42997 -> 24137
$ws.Cells.Item(41, 12).Value = 24137  # L41: 42997 -> 24137
$ws.Cells.Item(41, 14).Value = -24917  # N41: -43777 -> -24917
$ws.Cells.Item(97, 8).Value = 33332.332  # H97: 37499.25 -> 33332.332
$ws.Cells.Item(97, 10).Value = 33332.332  # J97: 37499.25 -> 33332.332
$ws.Cells.Item(97, 12).Value = 33332.332  # L97: 37499.25 -> 33332.332
$ws.Cells.Item(97, 14).Value = -35314.332  # N97: -39481.25 -> -35314.332
$ws.Cells.Item(126, 8).Value = 5805.3335  # H126: 4611.1816 -> 5805.3335
$ws.Cells.Item(126, 9).Value = 4124.8335  # I126: 3378 -> 4124.8335
$ws.Cells.Item(126, 10).Value = 9166.333000000001  # J126: 7899.6665 -> 9166.333000000001
$ws.Cells.Item(126, 11).Value = 12374.5005  # K126: 10134 -> 12374.5005
$ws.Cells.Item(126, 12).Value = 27498.999  # L126: 23698.9995 -> 27498.999
$ws.Cells.Item(126, 13).Value = -9904.500499999998  # M126: -7664 -> -9904.500499999998
$ws.Cells.Item(126, 14).Value = -32438.999  # N126: -28638.9995 -> -32438.999
$ws.Cells.Item(136, 8).Value = 1892.2667  # H136: 2065.1155 -> 1892.2667
$ws.Cells.Item(136, 9).Value = 1768.2307  # I136: 1982.4348 -> 1768.2307
$ws.Cells.Item(136, 10).Value = 2698.5  # J136: 2699 -> 2698.5
$ws.Cells.Item(136, 11).Value = 5304.6921  # K136: 5947.3044 -> 5304.6921
$ws.Cells.Item(136, 12).Value = 8095.5  # L136: 8097 -> 8095.5
$ws.Cells.Item(136, 13).Value = -2754.6921  # M136: -3397.3044 -> -2754.6921
$ws.Cells.Item(136, 14).Value = -13195.5  # N136: -13197 -> -13195.5
